$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.125.93"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "2.917.24"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'355.40"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'113.44"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "'39.68"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("E11").Value = "  +5.11%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "'20.05"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "'7.76"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "3.375.26"
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").Value = "2.903.77"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "'0.988"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").Value = "52.191.62"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'3.31"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "'7.60"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'14.10"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("D23").Value = "'71.28"
$ws.Range("E23").Value = "  +1.33%  "
$ws.Range("D24").Value = "'270.01"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").Value = "'0.181"
$ws.Range("E26").Value = "  +11.83%  "
$ws.Range("D27").Value = "'26.87"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'6.96"
$ws.Range("E29").Value = "  +13.93%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'10.66"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("D31").Value = "'0.104"
$ws.Range("E31").Value = "  +13.88%  "
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("E34").Value = "  +8.37%  "
$ws.Range("D35").Value = "'53.21"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").Value = "'0.0454"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  +6.31%  "
$ws.Range("D39").Value = "'18.85"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "'2.72"
$ws.Range("E41").Value = "  +8.03%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "'23.04"
$ws.Range("E43").Value = "  +4.52%  "
$ws.Range("E44").Value = "  -1.99%  "
$ws.Range("D45").Value = "'117.71"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "2.187.90"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("E49").Value = "  +14.63%  "
$ws.Range("D50").Value = "'0.0353"
$ws.Range("E50").Value = "  +12.44%  "
$ws.Range("D51").Value = "'0.956"
$ws.Range("E51").Value = "  -2.28%  "
